$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was collected for this product/vendor. Insert a
# fresh row right after the header block of data (row 9), pushing the
# existing rows 9-22 down to 10-23, then populate the new row with the
# latest week's figures (same vendor/region/category metadata, new date +
# prices).
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C9").Value = "Arica y Parinacota"
$ws.Range("D9").Value = 44498
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = 100112052
$ws.Range("G9").Value = "Albahaca"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 270
$ws.Range("K9").Value = 2000
$ws.Range("L9").Value = 2300
$ws.Range("M9").Value = 2150
$ws.Range("N9").Value = "`$/paquete"
$ws.Range("O9").Value = "Región de Arica y Parinacota"
$ws.Range("P9").Value = 2150
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = "Hortaliza"
